# Rename metadata sheet 'General' to 'Table' and update the active
# sheet/selection state to match the saved workbook view.
$wb = $excel.ActiveWorkbook

# Rename the "General" sheet to "Table".
$generalSheet = $wb.Worksheets.Item("General")
$generalSheet.Name = "Table"

# "Variables" sheet is no longer the active tab; its selection moved to K83.
$variablesSheet = $wb.Worksheets.Item("Variables")
$variablesSheet.Range("K83").Select()

# "Table" (formerly "General") becomes the active/selected sheet, with
# cell B1 selected.
$generalSheet.Activate()
$generalSheet.Range("B1").Select()
